$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-1651168788503177"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687908473377"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687908483384"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168790901999"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687909952497"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168788470268.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687884864688.csv"
$ws1.Range("B4").Value = "go_stims-1651168788488522.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687885020905.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511687906133595.csv"
$ws2.Range("B3").Value = "OB-16511687895888684.csv"
$ws2.Range("B4").Value = "ZB-match_1-16511687889557023.csv"
$ws2.Range("B5").Value = "TB-1651168790827463.csv"
$ws2.Range("B6").Value = "OB-16511687894708748.csv"
$ws2.Range("B7").Value = "ZB-match_0-16511687887473655.csv"
$ws2.Range("B8").Value = "ZB-match_6-16511687890541265.csv"
$ws2.Range("B9").Value = "TB-16511687902622168.csv"
$ws2.Range("B10").Value = "OB-16511687896532602.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651168790863413.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687908504066.csv"
$ws4.Range("B4").Value = "MM_stims-16511687908790789.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687908644133.csv"
$ws4.Range("B6").Value = "MM_stims-16511687909010339.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687908800852.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687909327872.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651168790948446.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687909793124.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687909067614.csv"
